$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at the front (A, B). This shifts the existing
# columns A-H (Power/Accuracy/.../Smart Sensors) to C-J, preserving
# their styles/formatting.
$ws.Range("A1:B1").EntireColumn.Insert()

# Pick up the header style (bold, centered, bordered) from the
# neighbouring header cell so the new header cells match the rest of
# row 1.
$ws.Range("C1").Copy($ws.Range("A1"))
$ws.Range("C1").Copy($ws.Range("B1"))

# New column headers for the inserted columns.
$ws.Range("A1").Value = "Pileup"
$ws.Range("B1").Value = "L1T Reduction Ratio"

# Populate rows 2-10 with the updated data set.
$ws.Range("A2").Value = 60
$ws.Range("B2").Value = 400
$ws.Range("C2").Value = 322345.06307234842097386718
$ws.Range("D2").Value = 0.25874125874125880609
$ws.Range("E2").Value = 0.25951903807615228637
$ws.Range("F2").Value = 0.25912956478239118630
$ws.Range("G2").Value = 0.00080469262310046703
$ws.Range("H2").Value = $false
$ws.Range("I2").Value = $false
$ws.Range("J2").Value = $false

$ws.Range("A3").Value = 200
$ws.Range("B3").Value = 53
$ws.Range("C3").Value = 51614519.03161904215812683105
$ws.Range("D3").Value = 0.44950706101785242730
$ws.Range("E3").Value = 0.44974673420421218406
$ws.Range("F3").Value = 0.44962686567164178442
$ws.Range("G3").Value = 0.00006538662603179312
$ws.Range("H3").Value = $false
$ws.Range("I3").Value = $false
$ws.Range("J3").Value = $false

$ws.Range("A4").Value = 200
$ws.Range("B4").Value = 53
$ws.Range("C4").Value = 25856177.96575045958161354065
$ws.Range("D4").Value = 0.45097255528910207278
$ws.Range("E4").Value = 0.45121300986403617417
$ws.Range("F4").Value = 0.45109275053304898062
$ws.Range("G4").Value = 0.00013095137997526509
$ws.Range("H4").Value = $true
$ws.Range("I4").Value = $false
$ws.Range("J4").Value = $false

$ws.Range("A5").Value = 200
$ws.Range("B5").Value = 53
$ws.Range("C5").Value = 51706478.77635879069566726685
$ws.Range("D5").Value = 0.83856382978723409405
$ws.Range("E5").Value = 0.84057584644094907045
$ws.Range("F5").Value = 0.83956863267208092783
$ws.Range("G5").Value = 0.00012210377243055911
$ws.Range("H5").Value = $false
$ws.Range("I5").Value = $true
$ws.Range("J5").Value = $false

$ws.Range("A6").Value = 200
$ws.Range("B6").Value = 53
$ws.Range("C6").Value = 40846779.62429585307836532593
$ws.Range("D6").Value = 0.44950706101785242730
$ws.Range("E6").Value = 0.44974673420421218406
$ws.Range("F6").Value = 0.44962686567164178442
$ws.Range("G6").Value = 0.00008262338634216193
$ws.Range("H6").Value = $false
$ws.Range("I6").Value = $false
$ws.Range("J6").Value = $true

$ws.Range("A7").Value = 200
$ws.Range("B7").Value = 53
$ws.Range("C7").Value = 25846604.50806888937950134277
$ws.Range("D7").Value = 0.83939757430361183754
$ws.Range("E7").Value = 0.83950946414289517872
$ws.Range("F7").Value = 0.83945351549483493425
$ws.Range("G7").Value = 0.00024368460951191799
$ws.Range("H7").Value = $true
$ws.Range("I7").Value = $true
$ws.Range("J7").Value = $false

$ws.Range("A8").Value = 200
$ws.Range("B8").Value = 53
$ws.Range("C8").Value = 40846779.62429585307836532593
$ws.Range("D8").Value = 0.44950706101785242730
$ws.Range("E8").Value = 0.44974673420421218406
$ws.Range("F8").Value = 0.44962686567164178442
$ws.Range("G8").Value = 0.00008262338634216193
$ws.Range("H8").Value = $false
$ws.Range("I8").Value = $true
$ws.Range("J8").Value = $true

$ws.Range("A9").Value = 200
$ws.Range("B9").Value = 53
$ws.Range("C9").Value = 20467713.46182410046458244324
$ws.Range("D9").Value = 0.45097255528910207278
$ws.Range("E9").Value = 0.45121300986403617417
$ws.Range("F9").Value = 0.45109275053304898062
$ws.Range("G9").Value = 0.00016542649924313101
$ws.Range("H9").Value = $true
$ws.Range("I9").Value = $false
$ws.Range("J9").Value = $true

$ws.Range("A10").Value = 200
$ws.Range("B10").Value = 53
$ws.Range("C10").Value = 20460139.28041300922632217407
$ws.Range("D10").Value = 0.83939757430361183754
$ws.Range("E10").Value = 0.83950946414289517872
$ws.Range("F10").Value = 0.83945351549483493425
$ws.Range("G10").Value = 0.00030783855576131761
$ws.Range("H10").Value = $true
$ws.Range("I10").Value = $true
$ws.Range("J10").Value = $true

